# Generate Report for Handback
# Rename the "zh-cn" locale sheet/table/column to "ko-kr" everywhere it
# appears, and refresh the handoff/handback timestamps that were recorded
# for this generated run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: rename the locale column header ("zh-cn" -> "ko-kr") ---
# The column name lives in the table's ListColumns metadata; in this
# engine that metadata only updates when the header cell text is edited
# while the table's header row is actually shown, so toggle ShowHeaders on
# for the edit and then restore the table back to its original (headerless)
# shape.
$overview = $wb.Worksheets.Item("Overview")
$loOverview = $overview.ListObjects.Item(1)
$loOverview.ShowHeaders = $true
$overview.Range("B1").Value = "ko-kr"
$loOverview.ShowHeaders = $false
$loOverview.Resize($overview.Range("A1:C1"))

# --- Detail sheet: refresh the handoff / handback timestamps ---
$detail = $wb.Worksheets.Item("zh-cn")
$detail.Range("E2:E5").Value = "2016-03-11 01:03:30"
$detail.Range("H2:H5").Value = "2016-03-19 01:30:54"

# --- Detail sheet: rename its table ("zh-cn"/"zh_cn" -> "ko-kr"/"ko_kr") ---
$loDetail = $detail.ListObjects.Item(1)
$loDetail.Name = "ko-kr"
$loDetail.DisplayName = "ko_kr"

# --- Rename the detail sheet itself last ---
$detail.Name = "ko-kr"

Write-Output "Renamed zh-cn to ko-kr and refreshed handback timestamps"
